$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 47, shifting existing rows 47:81 down to 48:82.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new weekly price record.
$ws.Cells.Item(47, 1).Value  = 2
$ws.Cells.Item(47, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(47, 3).Value  = "Coquimbo"
$ws.Cells.Item(47, 4).Value  = 44902
$ws.Cells.Item(47, 4).Style  = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat
$ws.Cells.Item(47, 5).Value  = 4
$ws.Cells.Item(47, 6).Value  = 100112026
$ws.Cells.Item(47, 7).Value  = "Haba"
$ws.Cells.Item(47, 8).Value  = "Sin especificar"
$ws.Cells.Item(47, 9).Value  = "Primera"
$ws.Cells.Item(47, 10).Value = 500
$ws.Cells.Item(47, 11).Value = 5000
$ws.Cells.Item(47, 12).Value = 6000
$ws.Cells.Item(47, 13).Value = 5500
$ws.Cells.Item(47, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(47, 16).Value = 220
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
